$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing hyperlink target addresses (column C, rows 2-6) before
# the column shift, in row order, so they can be re-attached to the correct
# cells (column B) after column A is removed.
$hlAddresses = @()
foreach ($hl in $ws.Hyperlinks) {
    $hlAddresses += $hl.Address
}

# Delete column A outright; B:G shift left to become A:F.
$ws.Range("A:A").Delete()

# The old hyperlinks still point at column C (their Range did not follow the
# shift), so drop them and recreate them on the correct cells (now column B).
$ws.Hyperlinks.Delete()

for ($i = 0; $i -lt $hlAddresses.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $ws.Hyperlinks.Add($cell, $hlAddresses[$i])
    # Re-apply the Hyperlink cell style so it matches the original formatting
    # (Hyperlinks.Add can otherwise leave a slightly different style applied).
    $cell.Style = "Hyperlink"
}
